$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Enemy Types")

# --- Fix K3 (attack speed multiplier) and AD3:AF3 (big-hit multipliers) ---
$ws.Range("K3").Value = 1.8
$ws.Range("AD3").Value = 5
$ws.Range("AE3").Value = 3
$ws.Range("AF3").Value = 0.8

# --- Fix the "damage / wait-time" formulas: they should multiply, not divide ---
$ws.Range("G4").Formula = '=F4*$G$3'
$ws.Range("G5:G23").Formula = '=F5*$G$3'

$ws.Range("L4").Formula = '=K4*L$3'
$ws.Range("L5:L23").Formula = '=K5*L$3'

$ws.Range("Q4").Formula = '=P4*Q$3'
$ws.Range("Q5:Q23").Formula = '=P5*Q$3'

$ws.Range("V4").Formula = '=U4*V$3'
$ws.Range("V5:V23").Formula = '=U5*V$3'

$ws.Range("AA4").Formula = '=Z4*AA$3'
$ws.Range("AA5:AA23").Formula = '=Z5*AA$3'

$ws.Range("AF4").Formula = '=AE4*AF$3'
$ws.Range("AF5:AF23").Formula = '=AE5*AF$3'

# --- New "Calculator DPS in WaitTime" block (rows 26-30) ---
# Write the labels in this specific order so the shared-string table
# gets populated in the same order as the target workbook.
$ws.Range("A26").Value = "Calculator DPS in WaitTime"
$ws.Range("A29").Value = "DPS"
$ws.Range("A28").Value = "Damage"
$ws.Range("A30").Value = "WaitTime"
$ws.Range("A27").Value = "HitSpeed"

$ws.Range("B27").Value = 0.5
$ws.Range("B28").Value = 10
$ws.Range("B29").Formula = "=B28*B27"
$ws.Range("B30").Formula = "=B28/B29"

$ws.Activate()
[void]$ws.Range("B27").Select()
